# Updated cryptos list: refresh Price (column D) and Volume(1h) (column E)
# values for each coin row. Values that look like plain numbers are
# prefixed with a leading apostrophe so Excel stores them as text
# (preserving formats like trailing zeros / multi-dot thousands values)
# instead of silently converting them to numeric cell values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "38.058.62"
$ws.Range("E2").Value = "  +2.42%  "
$ws.Range("D3").Value = "2.051.46"
$ws.Range("E3").Value = "  +1.52%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'229.07"
$ws.Range("E5").Value = "  +0.59%  "
$ws.Range("D6").Value = "'0.614"
$ws.Range("E6").Value = "  +1.53%  "
$ws.Range("D7").Value = "'60.79"
$ws.Range("E7").Value = "  +8.64%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "'0.385"
$ws.Range("E9").Value = "  +2.32%  "
$ws.Range("D10").Value = "'0.0803"
$ws.Range("E10").Value = "  +2.96%  "
$ws.Range("E11").Value = "  +1.71%  "
$ws.Range("D12").Value = "'14.77"
$ws.Range("E12").Value = "  +3.77%  "
$ws.Range("D13").Value = "2.352.79"
$ws.Range("E13").Value = "  +1.48%  "
$ws.Range("D14").Value = "'21.01"
$ws.Range("E14").Value = "  +5.20%  "
$ws.Range("D15").Value = "'5.33"
$ws.Range("E15").Value = "  +2.74%  "
$ws.Range("D16").Value = "'0.758"
$ws.Range("E16").Value = "  +2.76%  "
$ws.Range("D17").Value = "2.055.71"
$ws.Range("E17").Value = "  +2.12%  "
$ws.Range("D18").Value = "38.038.40"
$ws.Range("E18").Value = "  +2.63%  "
$ws.Range("D19").Value = "'6.31"
$ws.Range("E19").Value = "  +2.93%  "
$ws.Range("D20").Value = "'69.80"
$ws.Range("E20").Value = "  +1.30%  "
$ws.Range("D21").Value = "0.0₃0831"
$ws.Range("E21").Value = "  +1.87%  "
$ws.Range("D22").Value = "'226.07"
$ws.Range("E22").Value = "  +1.40%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").Value = "'2.43"
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("D25").Value = "'2.23"
$ws.Range("E25").Value = "  +1.57%  "
$ws.Range("D26").Value = "'165.96"
$ws.Range("E26").Value = "  +1.62%  "
$ws.Range("D27").Value = "'9.21"
$ws.Range("E27").Value = "  +2.30%  "
$ws.Range("D28").Value = "'0.133"
$ws.Range("E28").Value = "  +4.45%  "
$ws.Range("D29").Value = "'19.01"
$ws.Range("E29").Value = "  +1.70%  "
$ws.Range("E30").Value = "  -0.66%  "
$ws.Range("D31").Value = "'0.120"
$ws.Range("E31").Value = "  +1.93%  "
$ws.Range("D32").Value = "'4.53"
$ws.Range("E32").Value = "  +1.98%  "
$ws.Range("D33").Value = "'4.58"
$ws.Range("E33").Value = "  +2.96%  "
$ws.Range("D34").Value = "'2.04"
$ws.Range("E34").Value = "  +8.73%  "
$ws.Range("D35").Value = "'0.0605"
$ws.Range("E35").Value = "  +1.19%  "
$ws.Range("D36").Value = "'6.31"
$ws.Range("E36").Value = "  +15.41%  "
$ws.Range("E37").Value = "  -2.31%  "
$ws.Range("D38").Value = "'3.28"
$ws.Range("E38").Value = "  +3.50%  "
$ws.Range("E39").Value = "  +0.17%  "
$ws.Range("D40").Value = "1.519.13"
$ws.Range("E40").Value = "  +3.20%  "
$ws.Range("D41").Value = "'97.53"
$ws.Range("E41").Value = "  +3.46%  "
$ws.Range("D42").Value = "'16.94"
$ws.Range("E42").Value = "  +4.32%  "
$ws.Range("D43").Value = "'0.0216"
$ws.Range("E43").Value = "  +1.56%  "
$ws.Range("D44").Value = "'2.86"
$ws.Range("E44").Value = "  +2.62%  "
$ws.Range("D45").Value = "'0.0920"
$ws.Range("E45").Value = "  +1.02%  "
$ws.Range("D46").Value = "'1.13"
$ws.Range("E46").Value = "  +1.61%  "
$ws.Range("D47").Value = "'4.04"
$ws.Range("E47").Value = "  -4.36%  "
$ws.Range("E48").Value = "  +0.91%  "
$ws.Range("E49").Value = "  +1.22%  "
$ws.Range("D50").Value = "'7.02"
$ws.Range("E50").Value = "  -0.61%  "
$ws.Range("D51").Value = "2.241.82"
$ws.Range("E51").Value = "  +1.60%  "
